# Generate Report for Handback
#
# This localization-status report is refreshed after a handback event:
#  - the zh-cn / de-de "Publish" status moves from "In Translation" to
#    "Handed back: in sync with en-US"
#  - each language sheet grows two new facts per row: the "Latest Target
#    File" (a link back to the source .md) and the "Latest Handback File"
#    (the .xlf that came back from translation), plus a fresh "Latest
#    Handback DateTime"
#  - a few columns are widened so the new/longer values aren't truncated

$wb = $excel.ActiveWorkbook

$mdUrl9e4e = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b51aa9ebd0218bba4f04f8a2558ba61a94c60ffc/e2e/9e4e716d-4712-425f-b470-f38d1dc2d66a.md"
$mdUrlff31 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b51aa9ebd0218bba4f04f8a2558ba61a94c60ffc/e2e/ff31b10f-fc64-4802-b662-8a5d6469b576.md"
$name9e4e  = "9e4e716d-4712-425f-b470-f38d1dc2d66a.md"
$nameff31  = "ff31b10f-fc64-4802-b662-8a5d6469b576.md"

# ---------------------------------------------------------------------
# Overview sheet: handback is now in sync with en-US for both languages
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# Publish-status columns (E:F) get noticeably wider to fit the new text
$overview.Range("E:E").ColumnWidth = 29.166666666666668
$overview.Range("F:F").ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------
# zh-cn sheet: handback completed at 2016-08-20 16:35:35
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $mdUrl9e4e, [System.Type]::Missing, [System.Type]::Missing, $name9e4e)
$zhcn.Range("J2").Value = "9e4e716d-4712-425f-b470-f38d1dc2d66a.80bf98050e8108943142b2551f730c593e01a559.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-20 16:35:35"

$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $mdUrlff31, [System.Type]::Missing, [System.Type]::Missing, $nameff31)
$zhcn.Range("J3").Value = "ff31b10f-fc64-4802-b662-8a5d6469b576.96a592aa2c2594816308c06972ccfd4becc977cf.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-20 16:35:35"

$zhcn.Range("C:C").ColumnWidth = 29.166666666666668
$zhcn.Range("I:I").ColumnWidth = 39.166666666666664
$zhcn.Range("J:J").ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------
# de-de sheet: handback completed at 2016-08-20 16:35:41
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Hyperlinks.Add($dede.Range("I2"), $mdUrl9e4e, [System.Type]::Missing, [System.Type]::Missing, $name9e4e)
$dede.Range("J2").Value = "9e4e716d-4712-425f-b470-f38d1dc2d66a.80bf98050e8108943142b2551f730c593e01a559.de-de.xlf"
$dede.Range("K2").Value = "2016-08-20 16:35:41"

$dede.Hyperlinks.Add($dede.Range("I3"), $mdUrlff31, [System.Type]::Missing, [System.Type]::Missing, $nameff31)
$dede.Range("J3").Value = "ff31b10f-fc64-4802-b662-8a5d6469b576.96a592aa2c2594816308c06972ccfd4becc977cf.de-de.xlf"
$dede.Range("K3").Value = "2016-08-20 16:35:41"

$dede.Range("C:C").ColumnWidth = 29.166666666666668
$dede.Range("I:I").ColumnWidth = 39.166666666666664
$dede.Range("J:J").ColumnWidth = 39.166666666666664
